# Regenerate the handback status report with the latest run's identifiers.
#
# Two source documents were re-run by the localization pipeline and picked
# up new generated GUID-based file names + new hash + new timestamps:
#   doc1: 53478499-934e-480e-bd4b-871f57b73d9c  ->  90a9d022-6116-4a3b-bf78-763ea930d1ce
#   doc2: 839848ef-74c6-441f-b782-c0f541aa1a71  ->  ffffbf1238f6-29e7-4fe8-b526-03639834ce9e
# and the xliff content hash used in the generated handoff/handback file
# names changed from two distinct hashes to a single shared hash:
#   e2160e0be45c77815f671bc1b7d30101bbb330bc / 7432e153b21e08bcc6d0056b3d9978990a763004
#     -> 6539955a6181f4bcc191382f7e185d522e30c4aa

$wb = $excel.ActiveWorkbook

$oldDoc1 = "53478499-934e-480e-bd4b-871f57b73d9c"
$newDoc1 = "90a9d022-6116-4a3b-bf78-763ea930d1ce"
$oldDoc2 = "839848ef-74c6-441f-b782-c0f541aa1a71"
$newDoc2 = "ffffbf1238f6-29e7-4fe8-b526-03639834ce9e"

$newHash = "6539955a6181f4bcc191382f7e185d522e30c4aa"

$newXliffZhCn = $newDoc1 + "." + $newHash + ".zh-cn.xlf"
$newXliffDeDe = $newDoc1 + "." + $newHash + ".de-de.xlf"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = $newDoc1 + ".md"
$ws.Range("B2").Value = "e2e\" + $newDoc1 + ".md"
$ws.Range("G2").Value = "2016-09-05 13:17:36"

$ws.Range("A3").Value = $newDoc2 + ".md"
$ws.Range("B3").Value = "e2e\" + $newDoc2 + ".md"
$ws.Range("G3").Value = "2016-09-05 13:17:36"

$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3667e744e7d99241111e92340227cf2336251ff6/e2e/" + $oldDoc1 + ".md", "", "", "e2e\" + $newDoc1 + ".md")
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3667e744e7d99241111e92340227cf2336251ff6/e2e/" + $oldDoc2 + ".md", "", "", "e2e\" + $newDoc2 + ".md")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = $newDoc1 + ".md"
$ws.Range("G2").Value = $newXliffZhCn
$ws.Range("H2").Value = "2016-09-05 13:17:30"
$ws.Range("I2").Value = $newDoc1 + ".md"
$ws.Range("J2").Value = $newXliffZhCn
$ws.Range("K2").Value = "2016-09-05 13:17:57"

$ws.Range("A3").Value = $newDoc2 + ".md"
$ws.Range("G3").Value = $newXliffZhCn
$ws.Range("H3").Value = "2016-09-05 13:17:30"
$ws.Range("I3").Value = $newDoc2 + ".md"
$ws.Range("J3").Value = $newXliffZhCn
$ws.Range("K3").Value = "2016-09-05 13:17:57"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3667e744e7d99241111e92340227cf2336251ff6/e2e/" + $oldDoc1 + ".md", "", "", $newDoc1 + ".md")
$ws.Hyperlinks.Add($ws.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/18ba061699a28d479d5825d129da4b22170856c1/e2e/" + $oldDoc1 + ".md", "", "", $newDoc1 + ".md")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3667e744e7d99241111e92340227cf2336251ff6/e2e/" + $oldDoc2 + ".md", "", "", $newDoc2 + ".md")
$ws.Hyperlinks.Add($ws.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/18ba061699a28d479d5825d129da4b22170856c1/e2e/" + $oldDoc2 + ".md", "", "", $newDoc2 + ".md")

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = $newDoc1 + ".md"
$ws.Range("G2").Value = $newXliffDeDe
$ws.Range("H2").Value = "2016-09-05 13:17:36"
$ws.Range("I2").Value = $newDoc1 + ".md"
$ws.Range("J2").Value = $newXliffDeDe
$ws.Range("K2").Value = "2016-09-05 13:18:14"

$ws.Range("A3").Value = $newDoc2 + ".md"
$ws.Range("G3").Value = $newXliffDeDe
$ws.Range("H3").Value = "2016-09-05 13:17:36"
$ws.Range("I3").Value = $newDoc2 + ".md"
$ws.Range("J3").Value = $newXliffDeDe
$ws.Range("K3").Value = "2016-09-05 13:18:14"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3667e744e7d99241111e92340227cf2336251ff6/e2e/" + $oldDoc1 + ".md", "", "", $newDoc1 + ".md")
$ws.Hyperlinks.Add($ws.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/2d1d8848241ac10da3e0bf8c279e17ab03337b09/e2e/" + $oldDoc1 + ".md", "", "", $newDoc1 + ".md")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3667e744e7d99241111e92340227cf2336251ff6/e2e/" + $oldDoc2 + ".md", "", "", $newDoc2 + ".md")
$ws.Hyperlinks.Add($ws.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/2d1d8848241ac10da3e0bf8c279e17ab03337b09/e2e/" + $oldDoc2 + ".md", "", "", $newDoc2 + ".md")
